{"js": "// \"revisi bab v saran\" \u2014 replace the text of suggestion items 1 and 2 in\n// section 5.2 Saran (BAB V) with the revised wording from the diff, while\n// keeping the surrounding \"1. Purwarupa yang dibuat m\" / \"2. \" lead-ins\n// and paragraph formatting intact.\n\nconst body = context.document.body;\n\n// --- Suggestion item 1 ------------------------------------------------\n// Old: \"...masih m\" + \"asih menggunakan listrik sebagai sumber daya utama.\n//       Hal ini akan sangat merepotkan ketika tidak ada sumber listrik.\n//       Oleh karena itu diperlukan sumber daya utama selain listrik, salah\n//       satunya baterai supaya bisa mengisi kekurangan tersebut.\"\nconst oldItem1 =\n  \"asih menggunakan listrik sebagai sumber daya utama. Hal ini akan sangat \" +\n  \"merepotkan ketika tidak ada sumber listrik. Oleh karena itu diperlukan \" +\n  \"sumber daya utama selain listrik, salah satunya baterai supaya bisa \" +\n  \"mengisi kekurangan tersebut.\";\nconst newItem1 =\n  \"asih menggunakan sensor DHT11 untuk membaca suhu dan kelembaban. \" +\n  \"Terdapat sensor yang lebih akurat dibandingkan sensor DHT11 yaitu \" +\n  \"sensor DHT22 yang memiliki tingkat kesalahan sebesar 4% untuk \" +\n  \"pengukuran suhu dan 18% untuk pengukuran kelembaban. Sebaliknya sensor \" +\n  \"DHT11 memiliki tingkat kesalahan sebesar 1 - 7% untuk pengukuran suhu \" +\n  \"dan 11 - 35% untuk pengukuran kelembaban.\";\n\nconst found1 = body.search(oldItem1, { matchCase: true, matchWildcards: false });\nfound1.load(\"items\");\nawait context.sync();\n\nif (found1.items.length > 0) {\n  found1.items[0].insertText(newItem1, \"Replace\");\n  await context.sync();\n}\n\n// --- Suggestion item 2 ------------------------------------------------\n// Old: \"Hasil dari penelitian yang dilakukan masih dalam bentuk prototipe\n//       dan belum diterapkan ke kandang ayam yang sebenarnya. Di masa yang\n//       akan datang diharapkan hasil dari penelitian ini sudah bisa\n//       diterapkan dan digunakan ke kandang ayam yang sebenarnya.\"\nconst oldItem2 =\n  \"Hasil dari penelitian yang dilakukan masih dalam bentuk prototipe dan \" +\n  \"belum diterapkan ke kandang ayam yang sebenarnya. Di masa yang akan \" +\n  \"datang diharapkan hasil dari penelitian ini sudah bisa diterapkan dan \" +\n  \"digunakan ke kandang ayam yang sebenarnya.\";\nconst newItem2 =\n  \"Penelitian yang dilakukan menggunakan servo SG90 yang memiliki gir \" +\n  \"berbahan plastik serta memiliki kemampuan mengangkat beban seberat \\u00b1 \" +\n  \"1,8 kg. Terdapat servo yang lebih baik yaitu MG996 dengan gir berbahan \" +\n  \"metal dan mampu mengangkat beban \\u00b1 11 kg.\";\n\nconst found2 = body.search(oldItem2, { matchCase: true, matchWildcards: false });\nfound2.load(\"items\");\nawait context.sync();\n\nif (found2.items.length > 0) {\n  found2.items[0].insertText(newItem2, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# \"revisi bab v saran\" \u2014 replace the text of suggestion items 1 and 2 in\n# section 5.2 Saran (BAB V) with the revised wording from the diff, while\n# keeping the surrounding \"1. Purwarupa yang dibuat m\" / \"2. \" lead-ins\n# and paragraph formatting intact.\n\n$d = $word.ActiveDocument\n\n# --- Suggestion item 1 --------------------------------------------------\n# Old: \"...masih m\" + \"asih menggunakan listrik sebagai sumber daya utama.\n#       Hal ini akan sangat merepotkan ketika tidak ada sumber listrik.\n#       Oleh karena itu diperlukan sumber daya utama selain listrik, salah\n#       satunya baterai supaya bisa mengisi kekurangan tersebut.\"\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Text = \"asih menggunakan listrik sebagai sumber daya utama. Hal ini akan sangat merepotkan ketika tidak ada sumber listrik. Oleh karena itu diperlukan sumber daya utama selain listrik, salah satunya baterai supaya bisa mengisi kekurangan tersebut.\"\n$find1.Replacement.ClearFormatting()\n$find1.Replacement.Text = \"asih menggunakan sensor DHT11 untuk membaca suhu dan kelembaban. Terdapat sensor yang lebih akurat dibandingkan sensor DHT11 yaitu sensor DHT22 yang memiliki tingkat kesalahan sebesar 4% untuk pengukuran suhu dan 18% untuk pengukuran kelembaban. Sebaliknya sensor DHT11 memiliki tingkat kesalahan sebesar 1 - 7% untuk pengukuran suhu dan 11 - 35% untuk pengukuran kelembaban.\"\n$find1.Execute(\n    [ref]\"\",\n    [ref]$true,\n    [ref]$false,\n    [ref]$false,\n    [ref]$false,\n    [ref]$false,\n    [ref]$true,\n    [ref]1,\n    [ref]$false,\n    [ref]\"\",\n    [ref]2\n) | Out-Null\n\n# --- Suggestion item 2 --------------------------------------------------\n# Old: \"Hasil dari penelitian yang dilakukan masih dalam bentuk prototipe\n#       dan belum diterapkan ke kandang ayam yang sebenarnya. Di masa yang\n#       akan datang diharapkan hasil dari penelitian ini sudah bisa\n#       diterapkan dan digunakan ke kandang ayam yang sebenarnya.\"\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \"Hasil dari penelitian yang dilakukan masih dalam bentuk prototipe dan belum diterapkan ke kandang ayam yang sebenarnya. Di masa yang akan datang diharapkan hasil dari penelitian ini sudah bisa diterapkan dan digunakan ke kandang ayam yang sebenarnya.\"\n$find2.Replacement.ClearFormatting()\n$find2.Replacement.Text = \"Penelitian yang dilakukan menggunakan servo SG90 yang memiliki gir berbahan plastik serta memiliki kemampuan mengangkat beban seberat \u00b1 1,8 kg. Terdapat servo yang lebih baik yaitu MG996 dengan gir berbahan metal dan mampu mengangkat beban \u00b1 11 kg.\"\n$find2.Execute(\n    [ref]\"\",\n    [ref]$true,\n    [ref]$false,\n    [ref]$false,\n    [ref]$false,\n    [ref]$false,\n    [ref]$true,\n    [ref]1,\n    [ref]$false,\n    [ref]\"\",\n    [ref]2\n) | Out-Null\n"}
